$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.183.25"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").Value = "1.681.17"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("E4").Value = "  -0.48%  "

$ws.Range("D5").Value = "'211.71"
$ws.Range("E5").Value = "  -3.37%  "

$ws.Range("D6").Value = "'0.5302"
$ws.Range("E6").Value = "  -4.20%  "

$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("D8").Value = "'0.2690"
$ws.Range("E8").Value = "  -1.02%  "

$ws.Range("D9").Value = "'0.06326"
$ws.Range("E9").Value = "  -2.67%  "

$ws.Range("D10").Value = "'21.40"
$ws.Range("E10").Value = "  -3.22%  "

$ws.Range("D11").Value = "'0.07544"
$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.519"
$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.672.43"
$ws.Range("E13").Value = "  -1.05%  "

$ws.Range("D14").Value = "'0.5700"
$ws.Range("E14").Value = "  -2.59%  "

$ws.Range("D15").Value = "'0.000008145"
$ws.Range("E15").Value = "  -3.87%  "

$ws.Range("D16").Value = "'66.42"
$ws.Range("E16").Value = "  +1.31%  "

$ws.Range("D17").Value = "26.228.25"
$ws.Range("E17").Value = "  -1.14%  "

$ws.Range("D18").Value = "'4.876"
$ws.Range("E18").Value = "  -1.79%  "

$ws.Range("E19").Value = "  -0.48%  "

$ws.Range("D20").Value = "'10.59"
$ws.Range("E20").Value = "  -3.65%  "

$ws.Range("D21").Value = "'189.85"
$ws.Range("E21").Value = "  -0.57%  "

$ws.Range("D22").Value = "'6.229"
$ws.Range("E22").Value = "  -0.54%  "

$ws.Range("E23").Value = "  -0.42%  "

$ws.Range("D24").Value = "'149.08"
$ws.Range("E24").Value = "  -0.64%  "

$ws.Range("D25").Value = "'0.1264"
$ws.Range("E25").Value = "  -4.20%  "

$ws.Range("D26").Value = "'7.652"
$ws.Range("E26").Value = "  -3.92%  "

$ws.Range("D27").Value = "'15.90"
$ws.Range("E27").Value = "  +0.53%  "

$ws.Range("D28").Value = "'0.06471"
$ws.Range("E28").Value = "  +2.57%  "

$ws.Range("D29").Value = "'1.340"
$ws.Range("E29").Value = "  -5.07%  "

$ws.Range("D30").Value = "'1.289"
$ws.Range("E30").Value = "  -3.06%  "

$ws.Range("D31").Value = "'3.557"
$ws.Range("E31").Value = "  -0.98%  "

$ws.Range("D32").Value = "'3.561"
$ws.Range("E32").Value = "  -0.91%  "

$ws.Range("D33").Value = "'1.665"
$ws.Range("E33").Value = "  -0.82%  "

$ws.Range("D34").Value = "'1.011"
$ws.Range("E34").Value = "  -3.60%  "

$ws.Range("D35").Value = "'0.6114"
$ws.Range("E35").Value = "  -2.21%  "

$ws.Range("D36").Value = "'2.416"
$ws.Range("E36").Value = "  +0.41%  "

$ws.Range("D37").Value = "'2.717"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("D38").Value = "'6.210"
$ws.Range("E38").Value = "  -0.58%  "

$ws.Range("D39").Value = "'0.01614"
$ws.Range("E39").Value = "  -1.95%  "

$ws.Range("D40").Value = "1.101.23"
$ws.Range("E40").Value = "  -2.21%  "

$ws.Range("D41").Value = "'0.8717"
$ws.Range("E41").Value = "  -1.47%  "

$ws.Range("E42").Value = "  -0.94%  "

$ws.Range("D43").Value = "'100.09"
$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("D44").Value = "1.835.09"
$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("E45").Value = "  -1.89%  "

$ws.Range("D46").Value = "'57.21"
$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("D47").Value = "'1.006"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.05277"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.017"
$ws.Range("E49").Value = "  -2.56%  "

$ws.Range("D50").Value = "'0.4275"
$ws.Range("E50").Value = "  -0.59%  "

$ws.Range("D51").Value = "'5.984"
$ws.Range("E51").Value = "  -1.91%  "
